$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Logs" sheet: append a new test-mail row (row 30)
# ---------------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A30").Value = "Wil jij klant Dekker even bellen over die offerte?"
$logs.Range("B30").Value = "mailmind.test@zohomail.eu"
$logs.Range("C30").Value = "Testmail #20: Wil jij klant Dekker even bellen over die offerte?"
$logs.Range("D30").Value = "Samenwerking / Partnerverzoek"
$logs.Range("E30").Value = "Beste [Naam],`r`nBedankt voor je bericht. Ik zal klant Dekker zo snel mogelijk bellen over de offerte.`r`nMet vriendelijke groet,`r`n[Jouw Naam]"
$logs.Range("F30").Value = "2025-07-23 22:58:49"
$logs.Range("G30").Value = "Ja"
$logs.Range("H30").Value = "Nee"
$logs.Range("I30").Value = "Ja"
$logs.Range("J30").Value = "Nee"

# Re-fit the new row's height back to the sheet default (writing the
# multi-line "Antwoord" text otherwise leaves an explicit custom row height).
$logs.Rows.Item(30).AutoFit()

# Extend the conditional formatting ranges that highlight these columns so
# they keep covering the sheet's full data range (now through row 30).
foreach ($col in @("D", "G", "H", "I", "J")) {
    $oldRange = $col + "2:" + $col + "29"
    $newRange = $col + "2:" + $col + "30"
    $fcs = $logs.Range($oldRange).FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($logs.Range($newRange))
    }
}

# ---------------------------------------------------------------------------
# 2) "Dashboard" sheet: re-sync category counts
#    - rows 7-9 rotate to a new category order
#    - a new category/row (13) is appended
# ---------------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A7").Value = "Factuur / Administratie"
$dash.Range("A8").Value = "Offerte / Prijsaanvraag"
$dash.Range("A9").Value = "IT / Technisch probleem"

$dash.Range("A13").Value = "Samenwerking / Partnerverzoek"
$dash.Range("B13").Value = 1

# ---------------------------------------------------------------------------
# 3) Chart on the Dashboard sheet: extend the category/value series ranges
#    from row 12 to row 13 to include the newly added category.
# ---------------------------------------------------------------------------
$chartObj = $dash.ChartObjects().Item(1)
$chart = $chartObj.Chart
$ser = $chart.SeriesCollection(1)
$ser.Formula = "=SERIES('Dashboard'!`$B`$1,'Dashboard'!`$A`$2:`$A`$13,'Dashboard'!`$B`$2:`$B`$13,1)"
